$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1952.129
$ws.Range("I40").Value = 1249.125
$ws.Range("J40").Value = 2196.652
$ws.Range("K40").Value = 1249.125
$ws.Range("L40").Value = 2196.652
$ws.Range("M40").Value = -1074.125
$ws.Range("N40").Value = -2546.652

$ws.Range("H98").Value = 2348.875
$ws.Range("I98").Value = 2398.7144
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 2398.7144
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -900.7143999999998
$ws.Range("N98").Value = -4996

$ws.Range("H113").Value = 3228.0454
$ws.Range("I113").Value = 2129.6365
$ws.Range("J113").Value = 4326.4546
$ws.Range("K113").Value = 2129.6365
$ws.Range("L113").Value = 4326.4546
$ws.Range("M113").Value = 1124.3635
$ws.Range("N113").Value = -10834.4546

$ws.Range("H122").Value = 2348.875
$ws.Range("I122").Value = 2398.7144
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 7196.1432
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -4746.1432
$ws.Range("N122").Value = -10900

$ws.Range("H125").Value = 2500
$ws.Range("J125").Value = 2500
$ws.Range("L125").Value = 22500
$ws.Range("N125").Value = -27420

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2160.818
$ws.Range("I61").Value = 1296.258
$ws.Range("J61").Value = 4222.4614
$ws.Range("K61").Value = 1296.258
$ws.Range("L61").Value = 4222.4614
$ws.Range("M61").Value = -1084.258
$ws.Range("N61").Value = -4646.4614

$ws.Range("H74").Value = 2030.1875
$ws.Range("I74").Value = 1770.8182
$ws.Range("J74").Value = 2600.8
$ws.Range("K74").Value = 1770.8182
$ws.Range("L74").Value = 2600.8
$ws.Range("M74").Value = -896.8181999999999
$ws.Range("N74").Value = -4348.8

$ws.Range("H77").Value = 2030.1875
$ws.Range("I77").Value = 1770.8182
$ws.Range("J77").Value = 2600.8
$ws.Range("K77").Value = 8854.091
$ws.Range("L77").Value = 13004
$ws.Range("M77").Value = -4486.091
$ws.Range("N77").Value = -21740

$ws.Range("H108").Value = 32000
$ws.Range("J108").Value = 32000
$ws.Range("L108").Value = 32000
$ws.Range("N108").Value = -39680

$ws.Range("H136").Value = 2160.818
$ws.Range("I136").Value = 1296.258
$ws.Range("J136").Value = 4222.4614
$ws.Range("K136").Value = 3888.774
$ws.Range("L136").Value = 12667.3842
$ws.Range("M136").Value = -1338.774
$ws.Range("N136").Value = -17767.3842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 35000
$ws.Range("I62").Value = 20000
$ws.Range("K62").Value = 20000
$ws.Range("M62").Value = -19314

$ws.Range("H65").Value = 35000
$ws.Range("I65").Value = 20000
$ws.Range("K65").Value = 60000
$ws.Range("M65").Value = -56568

$ws.Range("H107").Value = 2244.0667
$ws.Range("I107").Value = 2332
$ws.Range("J107").Value = 1013
$ws.Range("K107").Value = 2332
$ws.Range("L107").Value = 1013
$ws.Range("M107").Value = -412
$ws.Range("N107").Value = -4853

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6946661.5
$ws.Range("I31").Value = 1570.4615
$ws.Range("J31").Value = 15154496
$ws.Range("K31").Value = 1570.4615
$ws.Range("L31").Value = 15154496
$ws.Range("M31").Value = -1275.4615
$ws.Range("N31").Value = -15155086

$ws.Range("H34").Value = 6946661.5
$ws.Range("I34").Value = 1570.4615
$ws.Range("J34").Value = 15154496
$ws.Range("K34").Value = 1570.4615
$ws.Range("L34").Value = 15154496
$ws.Range("M34").Value = -1368.4615
$ws.Range("N34").Value = -15154900

$ws.Range("H99").Value = 3619.3635
$ws.Range("I99").Value = 4266.6665
$ws.Range("J99").Value = 3376.625
$ws.Range("K99").Value = 4266.6665
$ws.Range("L99").Value = 3376.625
$ws.Range("M99").Value = -2768.6665
$ws.Range("N99").Value = -6372.625

$ws.Range("H109").Value = 25956.285
$ws.Range("J109").Value = 25956.285
$ws.Range("L109").Value = 25956.285
$ws.Range("N109").Value = -28036.285

$ws.Range("H126").Value = 3619.3635
$ws.Range("I126").Value = 4266.6665
$ws.Range("J126").Value = 3376.625
$ws.Range("K126").Value = 12799.9995
$ws.Range("L126").Value = 10129.875
$ws.Range("M126").Value = -10329.9995
$ws.Range("N126").Value = -15069.875

$ws.Range("H132").Value = 2453.0417
$ws.Range("I132").Value = 1720
$ws.Range("J132").Value = 3789.7646
$ws.Range("K132").Value = 5160
$ws.Range("L132").Value = 11369.2938
$ws.Range("M132").Value = -2630
$ws.Range("N132").Value = -16429.2938

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 2820.2
$ws.Range("I62").Value = 999
$ws.Range("J62").Value = 3022.5557
$ws.Range("K62").Value = 2997
$ws.Range("L62").Value = 9067.667099999999
$ws.Range("M62").Value = -2311
$ws.Range("N62").Value = -10439.6671

$ws.Range("H65").Value = 2820.2
$ws.Range("I65").Value = 999
$ws.Range("J65").Value = 3022.5557
$ws.Range("K65").Value = 8991
$ws.Range("L65").Value = 27203.0013
$ws.Range("M65").Value = -5559
$ws.Range("N65").Value = -34067.0013

$ws.Range("H68").Value = 885.4286
$ws.Range("I68").Value = 400
$ws.Range("J68").Value = 966.3333
$ws.Range("K68").Value = 1200
$ws.Range("L68").Value = 2898.9999
$ws.Range("M68").Value = -389
$ws.Range("N68").Value = -4520.9999

$ws.Range("H71").Value = 885.4286
$ws.Range("I71").Value = 400
$ws.Range("J71").Value = 966.3333
$ws.Range("K71").Value = 3600
$ws.Range("L71").Value = 8696.9997
$ws.Range("M71").Value = 456
$ws.Range("N71").Value = -16808.9997

$ws.Range("H80").Value = 2768.2
$ws.Range("J80").Value = 3425.7144
$ws.Range("L80").Value = 10277.1432
$ws.Range("N80").Value = -12149.1432

$ws.Range("H83").Value = 2768.2
$ws.Range("J83").Value = 3425.7144
$ws.Range("L83").Value = 30831.4296
$ws.Range("N83").Value = -40191.4296

$ws.Range("H97").Value = 1080.2667
$ws.Range("I97").Value = 1610
$ws.Range("J97").Value = 887.63635
$ws.Range("K97").Value = 4830
$ws.Range("L97").Value = 2662.90905
$ws.Range("M97").Value = -4334
$ws.Range("N97").Value = -3654.90905

$ws.Range("H98").Value = 2944.3333
$ws.Range("J98").Value = 4700
$ws.Range("L98").Value = 14100
$ws.Range("N98").Value = -17096

$ws.Range("H106").Value = 3342.7856
$ws.Range("J106").Value = 3342.7856
$ws.Range("L106").Value = 10028.3568
$ws.Range("N106").Value = -11920.3568

$ws.Range("H122").Value = 2620.775
$ws.Range("J122").Value = 3104.4688
$ws.Range("L122").Value = 27940.2192
$ws.Range("N122").Value = -32840.2192

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3700.3076
$ws.Range("I122").Value = 3332.4
$ws.Range("J122").Value = 4926.6665
$ws.Range("K122").Value = 9997.200000000001
$ws.Range("L122").Value = 14779.9995
$ws.Range("M122").Value = -7547.200000000001
$ws.Range("N122").Value = -19679.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 25349
$ws.Range("J109").Value = 25349
$ws.Range("L109").Value = 25349
$ws.Range("N109").Value = -28123
